$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.223.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.20"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6094"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.53%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07114"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2828"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07646"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.830.58"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.829"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6386"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001007"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.071.78"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.98%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.942"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.206.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.04%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.065"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.121"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1301"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.74"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06824"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.460"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.867"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.845"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.134"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6566"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.543"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.232.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.776"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01763"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.604"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9290"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.982.18"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.55"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000116"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.625"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.603"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.543"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1083"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.46%  "
